$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ between row 22 and row 23 and need to be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell22 = $ws.Range($col + "22")
    $cell23 = $ws.Range($col + "23")

    $tmp = $cell22.Value2
    $cell22.Value2 = $cell23.Value2
    $cell23.Value2 = $tmp
}
